# Apply updated dSF (column F) values for selected rows, per repull of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = 7
$ws.Range("F6").Value  = -5
$ws.Range("F9").Value  = -6
$ws.Range("F12").Value = -11
$ws.Range("F15").Value = -2
$ws.Range("F16").Value = -1
$ws.Range("F18").Value = 8
$ws.Range("F21").Value = 0
